$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 101, shifting existing rows 101-222 down to 102-223.
$ws.Rows(101).Insert()

# Populate the newly inserted row 101 with the new weekly data point.
$ws.Range("A101").Value = 10
$ws.Range("B101").Value = "Vega Modelo de Temuco"
$ws.Range("C101").Value = "La Araucanía"
$ws.Range("D101").Value = 44483
$ws.Range("E101").Value = 9
$ws.Range("F101").Value = "Fruta"
$ws.Range("G101").Value = 100108
$ws.Range("H101").Value = "Tropicales y subtropicales"
$ws.Range("I101").Value = 100108002
$ws.Range("J101").Value = "Mango"
$ws.Range("K101").Value = "Sin especificar"
$ws.Range("L101").Value = "Primera"
$ws.Range("M101").Value = 435
$ws.Range("N101").Value = 8000
$ws.Range("O101").Value = 9000
$ws.Range("P101").Value = 8575
$ws.Range("Q101").Value = "$/bandeja 4 kilos"
$ws.Range("R101").Value = "Perú"
$ws.Range("S101").Value = 2144
$ws.Range("T101").Value = 4
